$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @{
    11  = @("sd", "Statement-non-opinion")
    14  = @("sv", "Statement-opinion")
    27  = @("aa", "Agree/Accept")
    35  = @("sd", "Statement-non-opinion")
    40  = @("ba", "Appreciation")
    49  = @("sv", "Statement-opinion")
    60  = @("sd", "Statement-non-opinion")
    67  = @("sv", "Statement-opinion")
    71  = @("sd", "Statement-non-opinion")
    75  = @("sd", "Statement-non-opinion")
    79  = @("sd", "Statement-non-opinion")
    81  = @("ba", "Appreciation")
    103 = @("aa", "Agree/Accept")
    109 = @("sd", "Statement-non-opinion")
    110 = @("sd", "Statement-non-opinion")
    115 = @("aa", "Agree/Accept")
    139 = @("sd", "Statement-non-opinion")
    153 = @("sv", "Statement-opinion")
    167 = @("sd", "Statement-non-opinion")
    171 = @("%", "Uninterpretable")
    188 = @("sv", "Statement-opinion")
    210 = @("sv", "Statement-opinion")
    216 = @("aa", "Agree/Accept")
    222 = @("sv", "Statement-opinion")
    226 = @("sd", "Statement-non-opinion")
    231 = @("aa", "Agree/Accept")
    235 = @("sd", "Statement-non-opinion")
    251 = @("sv", "Statement-opinion")
    274 = @("sd", "Statement-non-opinion")
    320 = @("sd", "Statement-non-opinion")
}

foreach ($row in $changes.Keys) {
    $vals = $changes[$row]
    $ws.Cells.Item($row, 9).Value = $vals[0]
    $ws.Cells.Item($row, 10).Value = $vals[1]
}
